$wb = $excel.ActiveWorkbook

# --- Add "Errors" sheet after the last existing sheet, with the error message ---
$errSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$errSheet.Name = "Errors"
$errSheet.Range("A1").Value = "''Sheet ""Classes"" Row: 2 Invalid DDBNNN ""144Q1001""'"

# --- Add "Warnings" sheet after "Errors", becomes the active/selected sheet ---
$warnSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$warnSheet.Name = "Warnings"
$warnSheet.Range("D31").Select()
